$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update results for Steel (Iron & steel / Hydrogen value)
$ws.Range("B3").Value = 779.9074730793836

# Minor precision updates for Non-metallic minerals column
$ws.Range("D6").Value = 1891.100808345448
$ws.Range("D7").Value = 801.5679477565398
$ws.Range("D8").Value = 877.9100423599192
